$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item(1)

# Update F-column timestamps on the data sheet
$dataSheet.Range("F2").Value = "2021-10-05 14:21:24.789313"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:24.789321"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:24.789324"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:24.789326"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:24.789329"
$dataSheet.Range("F7").Value = "2021-10-05 14:21:24.789332"
$dataSheet.Range("F8").Value = "2021-10-05 14:21:24.789334"
$dataSheet.Range("F9").Value = "2021-10-05 14:21:24.789337"
$dataSheet.Range("F10").Value = "2021-10-05 14:21:24.789340"
$dataSheet.Range("F11").Value = "2021-10-05 14:21:24.789342"
$dataSheet.Range("F12").Value = "2021-10-05 14:21:24.789345"
$dataSheet.Range("F13").Value = "2021-10-05 14:21:24.789347"
$dataSheet.Range("F14").Value = "2021-10-05 14:21:24.789349"
$dataSheet.Range("F15").Value = "2021-10-05 14:21:24.789352"
$dataSheet.Range("F16").Value = "2021-10-05 14:21:24.789354"
$dataSheet.Range("F17").Value = "2021-10-05 14:21:24.789357"
$dataSheet.Range("F18").Value = "2021-10-05 14:21:24.789359"
$dataSheet.Range("F19").Value = "2021-10-05 14:21:24.789362"
$dataSheet.Range("F20").Value = "2021-10-05 14:21:24.789364"
$dataSheet.Range("F21").Value = "2021-10-05 14:21:24.789367"
$dataSheet.Range("F22").Value = "2021-10-05 14:21:24.789369"
$dataSheet.Range("F23").Value = "2021-10-05 14:21:24.789372"
$dataSheet.Range("F24").Value = "2021-10-05 14:21:24.789374"
$dataSheet.Range("F25").Value = "2021-10-05 14:21:24.789377"
$dataSheet.Range("F26").Value = "2021-10-05 14:21:24.789380"
$dataSheet.Range("F27").Value = "2021-10-05 14:21:24.789382"
$dataSheet.Range("F28").Value = "2021-10-05 14:21:24.789385"
$dataSheet.Range("F29").Value = "2021-10-05 14:21:24.789387"
$dataSheet.Range("F30").Value = "2021-10-05 14:21:24.789390"
$dataSheet.Range("F31").Value = "2021-10-05 14:21:24.789392"
$dataSheet.Range("F32").Value = "2021-10-05 14:21:24.789395"
$dataSheet.Range("F33").Value = "2021-10-05 14:21:24.789397"
$dataSheet.Range("F34").Value = "2021-10-05 14:21:24.789400"
$dataSheet.Range("F35").Value = "2021-10-05 14:21:24.789403"
$dataSheet.Range("F36").Value = "2021-10-05 14:21:24.789405"
$dataSheet.Range("F37").Value = "2021-10-05 14:21:24.789408"
$dataSheet.Range("F38").Value = "2021-10-05 14:21:24.789410"
$dataSheet.Range("F39").Value = "2021-10-05 14:21:24.789413"
$dataSheet.Range("F40").Value = "2021-10-05 14:21:24.789415"
$dataSheet.Range("F41").Value = "2021-10-05 14:21:24.789418"
$dataSheet.Range("F42").Value = "2021-10-05 14:21:24.789421"
$dataSheet.Range("F43").Value = "2021-10-05 14:21:24.789423"
$dataSheet.Range("F44").Value = "2021-10-05 14:21:24.789426"
$dataSheet.Range("F45").Value = "2021-10-05 14:21:24.789428"
$dataSheet.Range("F46").Value = "2021-10-05 14:21:24.789431"
$dataSheet.Range("F47").Value = "2021-10-05 14:21:24.789433"
$dataSheet.Range("F48").Value = "2021-10-05 14:21:24.789436"
$dataSheet.Range("F49").Value = "2021-10-05 14:21:24.789438"
$dataSheet.Range("F50").Value = "2021-10-05 14:21:24.789440"
$dataSheet.Range("F51").Value = "2021-10-05 14:21:24.789443"
$dataSheet.Range("F52").Value = "2021-10-05 14:21:24.789445"
$dataSheet.Range("F53").Value = "2021-10-05 14:21:24.789448"
$dataSheet.Range("F54").Value = "2021-10-05 14:21:24.789450"
$dataSheet.Range("F55").Value = "2021-10-05 14:21:24.789453"
$dataSheet.Range("F56").Value = "2021-10-05 14:21:24.789456"
$dataSheet.Range("F57").Value = "2021-10-05 14:21:24.789458"
$dataSheet.Range("F58").Value = "2021-10-05 14:21:24.789461"
$dataSheet.Range("F59").Value = "2021-10-05 14:21:24.789463"
$dataSheet.Range("F60").Value = "2021-10-05 14:21:24.789466"
$dataSheet.Range("F61").Value = "2021-10-05 14:21:24.789468"
$dataSheet.Range("F62").Value = "2021-10-05 14:21:24.789471"
$dataSheet.Range("F63").Value = "2021-10-05 14:21:24.789473"
$dataSheet.Range("F64").Value = "2021-10-05 14:21:24.789476"
$dataSheet.Range("F65").Value = "2021-10-05 14:21:24.789479"
$dataSheet.Range("F66").Value = "2021-10-05 14:21:24.789482"
$dataSheet.Range("F67").Value = "2021-10-05 14:21:24.789485"
$dataSheet.Range("F68").Value = "2021-10-05 14:21:24.789488"
$dataSheet.Range("F69").Value = "2021-10-05 14:21:24.789490"
$dataSheet.Range("F70").Value = "2021-10-05 14:21:24.789492"
$dataSheet.Range("F71").Value = "2021-10-05 14:21:24.789495"
$dataSheet.Range("F72").Value = "2021-10-05 14:21:24.789498"
$dataSheet.Range("F73").Value = "2021-10-05 14:21:24.789500"
$dataSheet.Range("F74").Value = "2021-10-05 14:21:24.789503"
$dataSheet.Range("F75").Value = "2021-10-05 14:21:24.789505"
$dataSheet.Range("F76").Value = "2021-10-05 14:21:24.789508"
$dataSheet.Range("F77").Value = "2021-10-05 14:21:24.789510"
$dataSheet.Range("F78").Value = "2021-10-05 14:21:24.789515"
$dataSheet.Range("F79").Value = "2021-10-05 14:21:24.789517"
$dataSheet.Range("F80").Value = "2021-10-05 14:21:24.789520"
$dataSheet.Range("F81").Value = "2021-10-05 14:21:24.789522"
$dataSheet.Range("F82").Value = "2021-10-05 14:21:24.789525"
$dataSheet.Range("F83").Value = "2021-10-05 14:21:24.789528"
$dataSheet.Range("F84").Value = "2021-10-05 14:21:24.789530"
$dataSheet.Range("F85").Value = "2021-10-05 14:21:24.789533"
$dataSheet.Range("F86").Value = "2021-10-05 14:21:24.789535"
$dataSheet.Range("F87").Value = "2021-10-05 14:21:24.789538"
$dataSheet.Range("F88").Value = "2021-10-05 14:21:24.789540"
$dataSheet.Range("F89").Value = "2021-10-05 14:21:24.789543"
$dataSheet.Range("F90").Value = "2021-10-05 14:21:24.789546"
$dataSheet.Range("F91").Value = "2021-10-05 14:21:24.789548"
$dataSheet.Range("F92").Value = "2021-10-05 14:21:24.789551"
$dataSheet.Range("F93").Value = "2021-10-05 14:21:24.789554"
$dataSheet.Range("F94").Value = "2021-10-05 14:21:24.789557"
$dataSheet.Range("F95").Value = "2021-10-05 14:21:24.789560"

# Add the metadata sheet after the data sheet
$newSheet = $wb.Worksheets.Add($null, $dataSheet)
$newSheet.Name = "metadata"

# Copy header-row and first-column styling from the data sheet so we reuse the existing style
$dataSheet.Range("B1:F1").Copy()
$newSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$newSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

# Data row
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Limb girdle muscular dystrophy"
$newSheet.Range("C2").Value = 185
$newSheet.Range("D2").Value = "'2.29"
$newSheet.Range("E2").Value = "2021-09-29T11:05:52.654397Z"
$newSheet.Range("F2").Value = "2021-10-05 14:21:24.786348"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/185/?format=json"

Write-Host "Edit complete"
